$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update existing date rows (row 2 and row 3 move forward one year / day) ---
$ws1.Range("A2").Value = 41309
$ws1.Range("A3").Value = 41310

# --- Fill in the workload hours for row 3 (B3/C3/D3) ---
# Set C3 first so "2h" becomes the first newly-added shared string, then B3/D3 add "1h".
$ws1.Range("C3").Value = "2h"
$ws1.Range("B3").Value = "1h"
$ws1.Range("D3").Value = "1h"

# Give B3 a centered (horizontal + vertical) look inside its thin box border.
$ws1.Range("B3").HorizontalAlignment = -4108
$ws1.Range("B3").VerticalAlignment = -4108

# Give C3/D3 a horizontally centered look inside their thin box border.
$ws1.Range("C3:D3").HorizontalAlignment = -4108

# --- Add two new day rows (4 and 5), copying the date cell's look from row 3 ---
$ws1.Range("A3").Copy($ws1.Range("A4"))
$ws1.Range("A3").Copy($ws1.Range("A5"))
$ws1.Range("A4").Value = 41311
$ws1.Range("A5").Value = 41312

# New rows' B/C/D cells get the same centered look (both axes) as B3
$ws1.Range("B3").Copy($ws1.Range("B4:D5"))
$ws1.Range("B4:D5").Value = ""

# --- Update the active selection to match the target workbook state ---
$ws1.Range("G6").Select()
